# Training Dashboard / Exam Dashboard refresh:
#  - re-run of the "days to expire" + "last update" computation (LAST UPDATE
#    moved from 08-Sep-2025 to 16-Sep-2025, i.e. +8 days, so PERIOD TO EXPIRE
#    drops by 8 for every training row; one SOP (row 16) rolls from VALID to
#    NOT VALID as a result)
#  - header / title font recolored to white-on-blue (single shared font)
#  - Exam Dashboard comments reworded + widened COMMENTS column

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write literal text into a cell without Excel's "looks like a date"
# auto-conversion clobbering it (and without disturbing the cell's existing
# style). Going through Formula -> Copy -> PasteSpecial(values) keeps the
# already-resolved text Variant instead of re-parsing a raw PowerShell string.
# ---------------------------------------------------------------------------
function Set-LiteralText {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)   # xlPasteValues
}

# ===========================================================================
# Sheet 1: Training Dashboard
# ===========================================================================
$ws1 = $wb.Worksheets.Item("Training Dashboard")

# --- Row 16 (Equipment Operation Procedure SOP-031) flips from VALID to
#     NOT VALID this cycle, so pick up the red "expired" row formatting used
#     by rows 15/17/18/19 first (copy formats only; values untouched here).
$ws1.Range("A17:K17").Copy()
$ws1.Range("A16:K16").PasteSpecial(-4122)   # xlPasteFormats

# --- PERIOD TO EXPIRE (col H) and LAST UPDATE (col I) refresh for every
#     data row (3-23): LAST UPDATE advances 8 days -> PERIOD TO EXPIRE -8.
for ($row = 3; $row -le 23; $row++) {
    $hCell = $ws1.Cells.Item($row, 8)
    $hCell.Value = $hCell.Value2 - 8

    Set-LiteralText $ws1.Cells.Item($row, 9) "16-Sep-2025"
}

# --- Row 16 now reads NOT VALID (period-to-expire went negative).
Set-LiteralText $ws1.Cells.Item(16, 10) "NOT VALID"

# --- Title (A1) and header row (row 2) both move to the same bold white
#     font (title loses its old 14pt size, header keeps its navy fill).
$title = $ws1.Range("A1")
$title.Font.Size = 11
$title.Font.Bold = $true
$title.Font.Color = RGB(255, 255, 255)

$header1 = $ws1.Range("A2:K2")
$header1.Font.Bold = $true
$header1.Font.Color = RGB(255, 255, 255)

# ===========================================================================
# Sheet 2: Exam Dashboard
# ===========================================================================
$ws2 = $wb.Worksheets.Item("Exam Dashboard")

# --- COMMENTS column (E) reworded from the terse "OK" to "date is valid",
#     and widened to fit.
for ($row = 3; $row -le 8; $row++) {
    Set-LiteralText $ws2.Cells.Item($row, 5) "date is valid"
}
$ws2.Columns.Item(5).ColumnWidth = 14.17

# --- Title (A1) and header row (row 2) get the same font refresh as sheet 1.
$title2 = $ws2.Range("A1")
$title2.Font.Size = 11
$title2.Font.Bold = $true
$title2.Font.Color = RGB(255, 255, 255)

$header2 = $ws2.Range("A2:G2")
$header2.Font.Bold = $true
$header2.Font.Color = RGB(255, 255, 255)

"done"
